$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing company (was "2" placeholder id) ---
$ws.Range("A2").Value = 'Bermuda'
$ws.Range("B2").Value = "'1"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = 'Brokerage & Investment Banking'
$ws.Range("D2").Value = 0.00369
$ws.Range("E2").Value = -0.22
$ws.Range("F2").Value = 0.12
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = -0.0203691866454293
$ws.Range("J2").Value = -0.01450471814146774
$ws.Range("K2").Value = 289.5
$ws.Range("L2").Value = 0.1198360791456246
$ws.Range("M2").Value = 356.6
$ws.Range("N2").Value = 0.08026289135474579
$ws.Range("O2").Value = 1.231778929188256
$ws.Range("P2").Value = 197
$ws.Range("Q2").Value = 0.04434040829188143
$ws.Range("R2").Value = 0.6804835924006909
$ws.Range("S2").Value = 159.6
$ws.Range("T2").Value = 0.4475602916432979
$ws.Range("U2").Value = 1109.7
$ws.Range("V2").Value = 0.2497692948299534
$ws.Range("W2").Value = 0.4719595696119987
$ws.Range("X2").Value = 0.04624340514724566
$ws.Range("Y2").Value = 0.425716164464753
$ws.Range("Z2").Value = 1.280885220620396
$ws.Range("AA2").Value = -0.01857887909667058
$ws.Range("AB2").Value = 0.03892414395907518
$ws.Range("AC2").Value = -0.05750302305574576
$ws.Range("AD2").Value = 1967.4
$ws.Range("AE2").Value = 691.5394054901406
$ws.Range("AF2").Value = 2658.93940549014
$ws.Range("AG2").Value = 1549.23940549014
$ws.Range("AH2").Value = 0.3744015111683072
$ws.Range("AI2").Value = 0.7734782773445161
$ws.Range("AJ2").Value = 0.2585452875263033
$ws.Range("AK2").Value = 0.6654981662479968
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 22.08080808080808
$ws.Range("AP2").Value = 17.38764764859866

# --- Row 3: replace LOM Financial Limited with Lazard Ltd data ---
$ws.Range("A3").Value = 'Bermuda'
$ws.Range("B3").Value = 'Lazard Ltd (NYSE:LAZ)'
$ws.Range("C3").Value = 'Brokerage & Investment Banking'
$ws.Range("D3").Value = 0.00369
$ws.Range("E3").Value = -0.22
$ws.Range("F3").Value = 0.12
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = -0.0203691866454293
$ws.Range("J3").Value = -0.01450471814146774
$ws.Range("K3").Value = 289.5
$ws.Range("L3").Value = 0.1198360791456246
$ws.Range("M3").Value = 356.6
$ws.Range("N3").Value = 0.08026289135474579
$ws.Range("O3").Value = 1.231778929188256
$ws.Range("P3").Value = 197
$ws.Range("Q3").Value = 0.04434040829188143
$ws.Range("R3").Value = 0.6804835924006909
$ws.Range("S3").Value = 159.6
$ws.Range("T3").Value = 0.4475602916432979
$ws.Range("U3").Value = 1109.7
$ws.Range("V3").Value = 0.2497692948299534
$ws.Range("W3").Value = 0.4719595696119987
$ws.Range("X3").Value = 0.04624340514724566
$ws.Range("Y3").Value = 0.425716164464753
$ws.Range("Z3").Value = 1.280885220620396
$ws.Range("AA3").Value = -0.01857887909667058
$ws.Range("AB3").Value = 0.03892414395907518
$ws.Range("AC3").Value = -0.05750302305574576
$ws.Range("AD3").Value = 1967.4
$ws.Range("AE3").Value = 691.5394054901406
$ws.Range("AF3").Value = 2658.93940549014
$ws.Range("AG3").Value = 1549.23940549014
$ws.Range("AH3").Value = 0.3744015111683072
$ws.Range("AI3").Value = 0.7734782773445161
$ws.Range("AJ3").Value = 0.2585452875263033
$ws.Range("AK3").Value = 0.6654981662479968
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = 0
$ws.Range("AN3").Value = 22.08080808080808
$ws.Range("AP3").Value = 17.38764764859866

# --- Row 4: remove entirely (Lazard Ltd row merged up / replaced) ---
$ws.Rows(4).Delete()
